$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 'maa://25390 (95.19), maa://24702 (94.79), maa://36681 (86.42)'
$ws.Range("L2").Value = 'maa://39402 (94.44), *maa://24633 (56.36), *maa://30515 (70.09), *maa://34787 (73.56), ***maa://20792 (11.93), ***maa://29083 (27.78)'
$ws.Range("T2").Value = 'maa://22742 (91.41), *maa://20791 (62.2)'
$ws.Range("AB2").Value = 'maa://21246 (91.45), maa://36684 (93.42), ***maa://22731 (6.25)'
$ws.Range("D3").Value = 'maa://40192 (96.84), maa://36987 (96.15), maa://39849 (88.89)'
$ws.Range("H3").Value = 'maa://21247 (98.41), *maa://22748 (60.0)'
$ws.Range("L3").Value = '*maa://22880 (65.33), maa://20276 (88.24), *maa://22749 (75.0)'
$ws.Range("P3").Value = 'maa://21249 (94.58), maa://26254 (97.14), **maa://22738 (50.0)'
$ws.Range("T3").Value = 'maa://24617 (90.62), maa://45854 (85.14), **maa://20790 (43.48), ***maa://37170 (16.67)'
$ws.Range("X3").Value = 'maa://27396 (84.15), maa://27484 (96.21), maa://27480 (83.78)'
$ws.Range("P4").Value = 'maa://49983 (94.83), maa://50121 (92.31)'
$ws.Range("T4").Value = 'maa://32509 (94.53), maa://27295 (87.5), maa://22754 (89.19), *maa://31008 (79.55), *maa://21746 (54.55)'
$ws.Range("X4").Value = '**maa://32495 (48.4), maa://43217 (93.04), ***maa://31785 (22.22), ***maa://36683 (29.79)'
$ws.Range("AB4").Value = '*maa://32658 (76.0)'
$ws.Range("AF4").Value = '*maa://30062 (66.07), *maa://39394 (66.67), ***maa://26209 (13.04)'
$ws.Range("D5").Value = 'maa://21245 (84.95), maa://22744 (81.48)'
$ws.Range("P6").Value = 'maa://31836 (94.12), maa://30381 (94.44)'
$ws.Range("T6").Value = 'maa://37411 (90.48)'
$ws.Range("D7").Value = 'maa://21955 (93.75)'
$ws.Range("P7").Value = 'maa://22750 (92.73)'
$ws.Range("T7").Value = 'maa://21291 (85.96)'
$ws.Range("X7").Value = 'maa://22399 (96.11), *maa://22758 (77.22)'
$ws.Range("A8").Value = '更新日期：2025.05.14 13:20:50'
$ws.Range("H8").Value = '*maa://24371 (55.13)'
$ws.Range("X8").Value = 'maa://21411 (95.93)'
$ws.Range("T9").Value = 'maa://26222 (98.28), **maa://22866 (30.19)'
$ws.Range("AB9").Value = 'maa://28711 (86.96), **maa://39938 (46.34), **maa://45044 (40.0), **maa://27377 (42.86), maa://40166 (94.59), ***maa://25174 (19.05)'
$ws.Range("AF9").Value = 'maa://26206 (88.55), *maa://22865 (51.85)'
$ws.Range("D10").Value = '***maa://25695 (19.0), ***maa://39951 (12.5), ***maa://34206 (22.22), *maa://45271 (59.7), ***maa://39243 (25.0), maa://54000 (100.0)'
$ws.Range("T10").Value = 'maa://27395 (96.57), maa://22755 (88.89), **maa://22756 (40.91), ***maa://21737 (11.76)'
$ws.Range("X10").Value = 'maa://22301 (97.87), maa://45828 (93.02), maa://22726 (100.0)'
$ws.Range("D11").Value = 'maa://36707 (99.36)'
$ws.Range("T11").Value = 'maa://22747 (90.8), maa://22501 (98.13), maa://45521 (90.91)'
$ws.Range("AB11").Value = 'maa://29912 (97.73), maa://22516 (87.36), *maa://20794 (55.26)'
$ws.Range("H12").Value = 'maa://21867 (90.81), **maa://45826 (33.33)'
$ws.Range("X12").Value = 'maa://22753 (91.71), *maa://21485 (75.68), maa://37962 (91.67)'
$ws.Range("D13").Value = 'maa://24999 (92.62), maa://36673 (92.94), maa://25001 (86.3)'
$ws.Range("X13").Value = 'maa://34957 (82.69), **maa://22768 (50.0)'
$ws.Range("AF13").Value = '**maa://22737 (38.22), maa://39883 (87.39), *maa://39885 (51.28)'
$ws.Range("L14").Value = 'maa://39841 (94.04), maa://26245 (96.57), maa://21288 (96.3), maa://36682 (95.74)'
$ws.Range("P14").Value = 'maa://23250 (98.82), maa://20107 (87.1), maa://22772 (100.0), *maa://22745 (66.67)'
$ws.Range("AB14").Value = 'maa://22764 (96.3)'
$ws.Range("D15").Value = '*maa://22743 (79.06), maa://22734 (84.43), *maa://30808 (65.22), *maa://36048 (66.0), maa://45058 (84.21)'
$ws.Range("L15").Value = '*maa://21334 (54.55)'
$ws.Range("P15").Value = 'maa://24762 (91.11), *maa://22727 (70.0)'
$ws.Range("T15").Value = 'maa://23892 (96.51)'
$ws.Range("AF15").Value = 'maa://21364 (81.27), maa://36666 (81.38), *maa://22766 (68.29)'
$ws.Range("T16").Value = 'maa://22729 (94.44), *maa://28648 (72.73), *maa://36674 (79.66)'
$ws.Range("AF16").Value = '*maa://23911 (68.38), maa://27755 (93.94)'
$ws.Range("T17").Value = '*maa://42324 (57.14)'
$ws.Range("D18").Value = 'maa://24570 (96.95)'
$ws.Range("H18").Value = 'maa://24421 (87.64)'
$ws.Range("L18").Value = 'maa://22466 (92.08), *maa://22732 (51.38), maa://52226 (92.86)'
$ws.Range("O18").NumberFormat = "@"
$ws.Range("O18").Value = '3'
$ws.Range("P18").Value = 'maa://24379 (100.0), maa://24380 (100.0), maa://54153 (100.0)'
$ws.Range("AB18").Value = 'maa://24393 (98.15)'
$ws.Range("D20").Value = 'maa://21432 (90.91), maa://25198 (94.44), *maa://20795 (50.76), maa://36680 (91.18)'
$ws.Range("L20").Value = 'maa://41331 (86.1)'
$ws.Range("T20").Value = 'maa://29113 (87.88)'
$ws.Range("X20").Value = 'maa://50085 (86.73), maa://49976 (86.25)'
$ws.Range("L21").Value = 'maa://31731 (96.55)'
$ws.Range("AB21").Value = 'maa://21443 (82.56), ***maa://23820 (30.0), **maa://52223 (41.18)'
$ws.Range("AF21").Value = 'maa://22524 (90.2), maa://22432 (81.9)'
$ws.Range("X22").Value = 'maa://21282 (98.7), *maa://37649 (71.05)'
$ws.Range("L23").Value = 'maa://39756 (95.45), maa://39875 (94.67)'
$ws.Range("X23").Value = '*maa://28503 (67.39)'
$ws.Range("D24").Value = '*maa://24368 (79.01), *maa://46650 (66.67)'
$ws.Range("X24").Value = 'maa://29988 (84.32), maa://23504 (93.58), **maa://22892 (41.29), *maa://25141 (77.37), *maa://36663 (77.78), ***maa://22815 (23.08), maa://52227 (100.0)'
$ws.Range("H25").Value = '*maa://29063 (73.12), *maa://25311 (74.11), ***maa://22725 (4.76), *maa://45047 (66.67)'
$ws.Range("P25").Value = 'maa://24382 (94.29)'
$ws.Range("AB25").Value = 'maa://31215 (88.81), *maa://24516 (79.35), maa://26001 (84.48)'
$ws.Range("H26").Value = 'maa://24913 (91.84)'
$ws.Range("AB26").Value = 'maa://42235 (95.68)'
$ws.Range("T27").Value = '*maa://30624 (76.06)'
$ws.Range("X28").Value = 'maa://39929 (91.88), maa://41749 (91.73), ***maa://39723 (13.89)'
$ws.Range("AF28").Value = 'maa://36660 (92.48), *maa://36701 (64.71)'
$ws.Range("H29").Value = '*maa://25175 (60.32)'
$ws.Range("L29").Value = 'maa://28432 (93.94), maa://28440 (83.46), maa://31400 (98.88), *maa://28650 (71.43)'
$ws.Range("O29").NumberFormat = "@"
$ws.Range("O29").Value = '3'
$ws.Range("P29").Value = '*maa://23168 (56.72), *maa://30050 (51.22), maa://54169 (100.0)'
$ws.Range("AF29").Value = '*maa://24080 (69.17), maa://42865 (81.55), ***maa://34960 (8.33)'
$ws.Range("T30").Value = '*maa://32940 (72.73), maa://24388 (94.74)'
$ws.Range("AB30").Value = 'maa://42979 (96.8), maa://45822 (100.0), maa://45045 (83.33)'
$ws.Range("H32").Value = 'maa://21895 (97.54), maa://36667 (98.02), **maa://20793 (38.0), maa://22760 (100.0)'
$ws.Range("T32").Value = 'maa://42859 (97.22), maa://41108 (88.0), maa://41238 (97.87), maa://45523 (100.0)'
$ws.Range("L35").Value = 'maa://41296 (97.14)'
$ws.Range("H37").Value = '*maa://24374 (57.14)'
$ws.Range("L37").Value = 'maa://45718 (98.31), *maa://47069 (80.0), maa://45789 (100.0)'
$ws.Range("P38").Value = '*maa://24383 (69.72)'
$ws.Range("T38").Value = 'maa://30713 (97.3)'
$ws.Range("AF38").Value = 'maa://36697 (88.85)'
$ws.Range("P41").Value = '**maa://35616 (40.0), maa://43177 (90.32)'
$ws.Range("P43").Value = '*maa://47403 (72.73)'
$ws.Range("T45").Value = '**maa://39364 (45.28)'
$ws.Range("H46").Value = 'maa://35931 (92.33), maa://43901 (94.34)'
$ws.Range("H47").Value = 'maa://27410 (96.84), maa://29661 (97.56), maa://28038 (84.62)'
$ws.Range("H52").Value = 'maa://24376 (96.97)'
$ws.Range("H53").Value = 'maa://32534 (94.83), **maa://32434 (33.33)'
$ws.Range("H59").Value = 'maa://31270 (94.44), maa://27746 (82.91)'
$ws.Range("H60").Value = '*maa://40438 (72.97)'
